$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Treatment query text (cell B5) needs an extra filter clause appended to
# its WHERE condition, matching the same pattern already used by the other
# sibling queries (e.g. Diagnoses / Treatment Response / Survival tabs).
$cell = $ws.Range("B5")
$oldText = $cell.Value()

$find = "WHERE `n    std.dbgap_accession = 'phs002431' AND prt.sex_at_birth = 'Male'`nORDER BY `n    trt.treatment_id ASC"
$replace = "WHERE `n    std.dbgap_accession = 'phs002431' AND prt.sex_at_birth = 'Male' AND trt.treatment_id  IS NOT NULL`nORDER BY `n    trt.treatment_id ASC"

$newText = $oldText.Replace($find, $replace)
$cell.Value = $newText

# Reflect where the user ended up working/saving: scrolled back to the top
# of the sheet with B2 selected (instead of the previous scrolled-down C7).
$ws.Range("B2").Select()
